$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.845.57"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "1.870.68"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'300.99"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'0.5321"
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("D8").Value = "'0.3748"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").Value = "'0.07178"
$ws.Range("D10").Value = "'21.49"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "'0.8857"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").Value = "'0.08155"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "1.843.23"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").Value = "'93.09"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "'5.257"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'14.68"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'0.000008526"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "26.829.79"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "'4.966"
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").Value = "'6.380"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("D24").Value = "'146.99"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").Value = "'2.252"
$ws.Range("E25").Value = "  -3.30%  "
$ws.Range("D26").Value = "'1.731"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").Value = "'18.03"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "'114.24"
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("D29").Value = "'4.735"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").Value = "'4.575"
$ws.Range("E30").Value = "  -6.54%  "
$ws.Range("D31").Value = "'0.09117"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").Value = "'0.7980"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("D33").Value = "'0.05001"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "'2.982"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").Value = "'1.168"
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("D36").Value = "'0.6006"
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("D37").Value = "'2.594"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").Value = "'3.154"
$ws.Range("E38").Value = "  -6.22%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.073"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01950"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").Value = "'6.616"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "'8.854"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").Value = "'115.37"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").Value = "'0.5121"
$ws.Range("E44").Value = "  +5.04%  "
$ws.Range("D45").Value = "'0.1495"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "'9.966"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("D48").Value = "'1.622"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").Value = "'37.58"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "'0.06015"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'62.10"
